$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.948.89'
$ws.Range('E2').Value = '  +0.15%  '
$ws.Range('D3').Value = '2.789.76'
$ws.Range('E3').Value = '  -1.38%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = "'358.43"
$ws.Range('E5').Value = '  +0.29%  '
$ws.Range('D6').Value = "'109.71"
$ws.Range('E6').Value = '  -2.36%  '
$ws.Range('E7').Value = '  -1.28%  '
$ws.Range('E9').Value = '  -1.88%  '
$ws.Range('D10').Value = "'40.36"
$ws.Range('E10').Value = '  -1.87%  '
$ws.Range('E11').Value = '  +1.91%  '
$ws.Range('E12').Value = '  -0.97%  '
$ws.Range('D13').Value = "'19.52"
$ws.Range('E13').Value = '  -2.85%  '
$ws.Range('E14').Value = '  -2.36%  '
$ws.Range('D15').Value = '3.225.41'
$ws.Range('E15').Value = '  -1.61%  '
$ws.Range('B16').Value = 'Polygon'
$ws.Range('C16').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D16').Value = "'0.950"
$ws.Range('E16').Value = '  +2.30%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '2.764.82'
$ws.Range('E17').Value = '  -2.23%  '
$ws.Range('D18').Value = '51.878.79'
$ws.Range('E18').Value = '  +0.14%  '
$ws.Range('D19').Value = "'7.44"
$ws.Range('E19').Value = '  -1.08%  '
$ws.Range('E20').Value = '  -2.33%  '
$ws.Range('D21').Value = "'13.20"
$ws.Range('E21').Value = '  -1.74%  '
$ws.Range('D22').Value = '0.0₃0977'
$ws.Range('E22').Value = '  -1.41%  '
$ws.Range('D23').Value = "'270.70"
$ws.Range('E23').Value = '  +0.45%  '
$ws.Range('D24').Value = "'70.11"
$ws.Range('E24').Value = '  +0.19%  '
$ws.Range('E25').Value = '  -1.95%  '
$ws.Range('D26').Value = "'26.48"
$ws.Range('E26').Value = '  -1.93%  '
$ws.Range('D28').Value = "'0.166"
$ws.Range('E28').Value = '  +18.77%  '
$ws.Range('D29').Value = "'10.27"
$ws.Range('E29').Value = '  -0.42%  '
$ws.Range('E30').Value = '  -4.91%  '
$ws.Range('D31').Value = "'52.08"
$ws.Range('E31').Value = '  -1.61%  '
$ws.Range('D32').Value = "'34.91"
$ws.Range('E32').Value = '  -1.87%  '
$ws.Range('E33').Value = '  -0.77%  '
$ws.Range('D34').Value = "'5.83"
$ws.Range('E34').Value = '  -0.96%  '
$ws.Range('E35').Value = '  +0.36%  '
$ws.Range('D36').Value = "'5.21"
$ws.Range('E36').Value = '  -4.24%  '
$ws.Range('E37').Value = '  -0.07%  '
$ws.Range('D38').Value = "'18.77"
$ws.Range('E38').Value = '  +0.82%  '
$ws.Range('D39').Value = "'3.20"
$ws.Range('E39').Value = '  -2.43%  '
$ws.Range('E40').Value = '  -3.46%  '
$ws.Range('E41').Value = '  +1.37%  '
$ws.Range('D42').Value = "'0.114"
$ws.Range('E42').Value = '  -1.71%  '
$ws.Range('D43').Value = "'2.23"
$ws.Range('E43').Value = '  -2.29%  '
$ws.Range('D44').Value = "'119.29"
$ws.Range('E44').Value = '  -3.74%  '
$ws.Range('D45').Value = "'21.85"
$ws.Range('E45').Value = '  -6.75%  '
$ws.Range('D46').Value = '2.077.85'
$ws.Range('E46').Value = '  -1.02%  '
$ws.Range('D47').Value = "'3.29"
$ws.Range('E47').Value = '  -2.49%  '
$ws.Range('E48').Value = '  -0.63%  '
$ws.Range('D49').Value = "'5.80"
$ws.Range('E49').Value = '  -2.75%  '
$ws.Range('D50').Value = "'0.950"
$ws.Range('E50').Value = '  -3.01%  '
$ws.Range('D51').Value = "'1.14"
$ws.Range('E51').Value = '  +31.59%  '
